$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 103: 72. Edit Distance ----
$ws.Range("A103").Value = "72. Edit Distance"
$ws.Range("B103").Value = "Medium"
$ws.Range("B103").Interior.Color = 49407
$ws.Range("C103").Value = "Dynamic Programming"
$ws.Range("D103").Value = "The indices represent the characters in the string up to that point,and the value represents the minimum operations to get there.The base case of 1 empty word is to either delete every character, or copy and add each character. Also, another base case is if the characters up to the point match, which will take 0 operations. The operations can be defined as index operations: matching (i+1, j+1), insert (i, j+1), delete (i+1, j), replace (i+1, j+1) after replacing. We represent it as a 2d dp grid and solve bottom-up (from bottom right to top left). The extra layer at the end for the base cases (empty strings)."
$ws.Hyperlinks.Add($ws.Range("E103"), "https://leetcode.com/problems/edit-distance/solutions/25849/java-dp-solution-o-nm/ ", [Type]::Missing, [Type]::Missing, "https://leetcode.com/problems/edit-distance/solutions/25849/java-dp-solution-o-nm/ ") | Out-Null

# ---- Row 104: 7. Reverse Integer ----
$ws.Range("A104").Value = "7. Reverse Integer"
$ws.Range("B104").Value = "Medium"
$ws.Range("B104").Interior.Color = 49407
$ws.Range("C104").Value = "Bit Manipulation"
$ws.Range("D104").Value = "To detect the overflow, we check if the resulting reversed integer is equal to the max integer except for the last digit (by chopping off the end by 10). Then we compare the last digit to see if it goes out of bounds. When an integer overflows, recognize that it is rolled over. Initialize the resulting num from 0, and build it up with checks before each step. A key implementation detail is to handle the negative before and after the core logic."
$ws.Hyperlinks.Add($ws.Range("E104"), "https://leetcode.com/problems/reverse-integer/solutions/3136892/java-beat-100-well-explained-code/ ", [Type]::Missing, [Type]::Missing, "https://leetcode.com/problems/reverse-integer/solutions/3136892/java-beat-100-well-explained-code/ ") | Out-Null

# ---- Row 105: 312. Burst Balloons ----
$ws.Range("A105").Value = "312. Burst Balloons"
$ws.Range("B105").Value = "Hard"
$ws.Range("B105").Interior.Color = 255
$ws.Range("C105").Value = "Dynamic Programming"
$ws.Hyperlinks.Add($ws.Range("E105"), "https://leetcode.com/problems/burst-balloons/solutions/76228/share-some-analysis-and-explanations/ ", [Type]::Missing, [Type]::Missing, "https://leetcode.com/problems/burst-balloons/solutions/76228/share-some-analysis-and-explanations/ ") | Out-Null
$ws.Range("D105").Value = "Review. The Brute Force solution is the decision tree O(N^n). The crux is to consider for the balloons to compute what happens if it is popped last, for the purpose of the dp array. The time complexity of the optimal is O(n^3) and O(n^2) space."

# ---- Fix hyperlink cell styling to reuse the existing "Hyperlink" cell style ----
# (Hyperlinks.Add creates a fresh style; re-apply formats copied from an already
# correctly-styled hyperlink cell so the three new cells share the original style.)
$ws.Range("E102").Copy() | Out-Null
$ws.Range("E103:E105").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- Resize the table / autofilter to cover the new rows ----
$tbl = $ws.ListObjects("Table2")
$tbl.Resize($ws.Range("A1:E105"))

$ws.Range("D111").Select()
